# Update observation counts (column E) and recompute shares (column F)
# F is defined as E / SUM(E) within each Gender+Period block:
#   rows 2-10, 11-19, 20-28, 29-37

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New counts for the rows whose N (E) value changed
$updatedE = @{
    2  = 9488
    3  = 18738
    5  = 3165
    8  = 2547
    9  = 1119
    13 = 9256
    14 = 2327
    15 = 3928
    16 = 1062
    20 = 11427
    21 = 30031
    23 = 4667
    24 = 4296
    25 = 1094
    26 = 3687
    29 = 5847
    30 = 22934
    31 = 9115
    32 = 3555
    35 = 4315
}

foreach ($row in $updatedE.Keys) {
    $ws.Range("E$row").Value = $updatedE[$row]
}

# Recompute P (column F) for every block using the updated N values
$blocks = @(
    ,@(2, 10)
    ,@(11, 19)
    ,@(20, 28)
    ,@(29, 37)
)

foreach ($block in $blocks) {
    $startRow = $block[0]
    $endRow = $block[1]

    $sum = 0
    for ($r = $startRow; $r -le $endRow; $r++) {
        $sum += $ws.Range("E$r").Value2
    }

    for ($r = $startRow; $r -le $endRow; $r++) {
        $n = $ws.Range("E$r").Value2
        $ws.Range("F$r").Value = $n / $sum
    }
}
